# Apply scheduled-runner price/profit updates to the Leve profit sheets.
# Each row's H-N columns (currentAveragePrice.. / LevePriceNQ/HQ / LeveProfitNQ/HQ)
# are refreshed from the latest market-board snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3001
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 3001
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 3001
$ws.Range("N51").Value = -3969
$ws.Range("H96").Value = 643.7273
$ws.Range("I96").Value = 947.8
$ws.Range("J96").Value = 390.33334
$ws.Range("K96").Value = 2843.4
$ws.Range("L96").Value = 1171.00002
$ws.Range("M96").Value = -1470.4
$ws.Range("N96").Value = -3917.00002
$ws.Range("H125").Value = 1444.3334
$ws.Range("I125").Value = 2163
$ws.Range("J125").Value = 1085
$ws.Range("K125").Value = 19467
$ws.Range("L125").Value = 9765
$ws.Range("M125").Value = -17007
$ws.Range("N125").Value = -14685
$ws.Range("H137").Value = 2876.492
$ws.Range("I137").Value = 2203.8928
$ws.Range("J137").Value = 3414.5715
$ws.Range("K137").Value = 6611.678400000001
$ws.Range("L137").Value = 10243.7145
$ws.Range("M137").Value = -4061.678400000001
$ws.Range("N137").Value = -15343.7145
$ws.Range("H138").Value = 3316.0923
$ws.Range("I138").Value = 3346.1738
$ws.Range("J138").Value = 3299.6191
$ws.Range("K138").Value = 10038.5214
$ws.Range("L138").Value = 9898.8573
$ws.Range("M138").Value = -4898.5214
$ws.Range("N138").Value = -20178.8573
$ws.Range("H141").Value = 4559.778
$ws.Range("I141").Value = 3141.682
$ws.Range("J141").Value = 10799.4
$ws.Range("K141").Value = 9425.045999999998
$ws.Range("L141").Value = 32398.2
$ws.Range("M141").Value = -4245.045999999998
$ws.Range("N141").Value = -42758.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1336.5
$ws.Range("I2").Value = 1393.0667
$ws.Range("J2").Value = 488
$ws.Range("K2").Value = 1393.0667
$ws.Range("L2").Value = 488
$ws.Range("M2").Value = -1280.0667
$ws.Range("N2").Value = -714
$ws.Range("H32").Value = 3492.25
$ws.Range("I32").Value = 1909.3478
$ws.Range("J32").Value = 39899
$ws.Range("K32").Value = 1909.3478
$ws.Range("L32").Value = 39899
$ws.Range("M32").Value = -1622.3478
$ws.Range("N32").Value = -40473
$ws.Range("H74").Value = 2389.3
$ws.Range("I74").Value = 2128.5715
$ws.Range("J74").Value = 2997.6667
$ws.Range("K74").Value = 2128.5715
$ws.Range("L74").Value = 2997.6667
$ws.Range("M74").Value = -1254.5715
$ws.Range("N74").Value = -4745.6667
$ws.Range("H77").Value = 2389.3
$ws.Range("I77").Value = 2128.5715
$ws.Range("J77").Value = 2997.6667
$ws.Range("K77").Value = 10642.8575
$ws.Range("L77").Value = 14988.3335
$ws.Range("M77").Value = -6274.8575
$ws.Range("N77").Value = -23724.3335
$ws.Range("H116").Value = 1336.5
$ws.Range("I116").Value = 1393.0667
$ws.Range("J116").Value = 488
$ws.Range("K116").Value = 1393.0667
$ws.Range("L116").Value = 488
$ws.Range("M116").Value = 900.9332999999999
$ws.Range("N116").Value = -5076
$ws.Range("H122").Value = 3150.4583
$ws.Range("I122").Value = 3363.5625
$ws.Range("J122").Value = 2724.25
$ws.Range("K122").Value = 10090.6875
$ws.Range("L122").Value = 8172.75
$ws.Range("M122").Value = -7640.6875
$ws.Range("N122").Value = -13072.75
$ws.Range("H132").Value = 4724.75
$ws.Range("I132").Value = 4299.6665
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 12898.9995
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -10368.9995
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1336.5
$ws.Range("I3").Value = 1393.0667
$ws.Range("J3").Value = 488
$ws.Range("K3").Value = 1393.0667
$ws.Range("L3").Value = 488
$ws.Range("M3").Value = -1279.0667
$ws.Range("N3").Value = -716
$ws.Range("H58").Value = 107022.86
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 107022.86
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 107022.86
$ws.Range("N58").Value = -107610.86
$ws.Range("M58").Value = $null
$ws.Range("H81").Value = 44895.668
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 44895.668
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 44895.668
$ws.Range("N81").Value = -47017.668
$ws.Range("H84").Value = 44895.668
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 44895.668
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 134687.004
$ws.Range("N84").Value = -145295.004
$ws.Range("H134").Value = 2968.8572
$ws.Range("I134").Value = 2968.8572
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 8906.571599999999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -6371.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 40000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 40000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 40000
$ws.Range("N74").Value = -41748
$ws.Range("H77").Value = 40000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 40000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 120000
$ws.Range("N77").Value = -128736
$ws.Range("H134").Value = 1894.72
$ws.Range("I134").Value = 1936.9048
$ws.Range("J134").Value = 1673.25
$ws.Range("K134").Value = 5810.7144
$ws.Range("L134").Value = 5019.75
$ws.Range("M134").Value = -3275.7144
$ws.Range("N134").Value = -10089.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7039526.5
$ws.Range("I131").Value = 1111111
$ws.Range("J131").Value = 7578473.5
$ws.Range("K131").Value = 3333333
$ws.Range("L131").Value = 22735420.5
$ws.Range("M131").Value = -3328293
$ws.Range("N131").Value = -22745500.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = $null
$ws.Range("H132").Value = 1507.32
$ws.Range("I132").Value = 1507.32
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4521.96
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1991.96

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 601.2727
$ws.Range("I55").Value = 642.1429000000001
$ws.Range("J55").Value = 529.75
$ws.Range("K55").Value = 642.1429000000001
$ws.Range("L55").Value = 529.75
$ws.Range("M55").Value = -469.1429000000001
$ws.Range("N55").Value = -875.75
$ws.Range("H132").Value = 2374.818
$ws.Range("I132").Value = 2309.25
$ws.Range("J132").Value = 2549.6667
$ws.Range("K132").Value = 6927.75
$ws.Range("L132").Value = 7649.000100000001
$ws.Range("M132").Value = -4397.75
$ws.Range("N132").Value = -12709.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 77942.836
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 77942.836
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 77942.836
$ws.Range("N46").Value = -78404.836
$ws.Range("H76").Value = 30000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 30000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30630
$ws.Range("H79").Value = 30000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 30000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32184
$ws.Range("H128").Value = 50715
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 50715
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 50715
$ws.Range("N128").Value = -60675
$ws.Range("H132").Value = 2241.25
$ws.Range("I132").Value = 2245.0908
$ws.Range("J132").Value = 2199
$ws.Range("K132").Value = 6735.2724
$ws.Range("L132").Value = 6597
$ws.Range("M132").Value = -4205.2724
$ws.Range("N132").Value = -11657
$ws.Range("H134").Value = 77942.836
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 77942.836
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 233828.508
$ws.Range("N134").Value = -238898.508
$ws.Range("H135").Value = 67450
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 67450
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 67450
$ws.Range("N135").Value = -77590

Write-Host "Updated 33 leve-profit rows across 8 sheets."
